# Weekly update: insert a new price observation as row 105 (Membrillo,
# Feria Lagunitas de Puerto Montt), pushing the existing rows 105-163
# down to 106-164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 105; Excel shifts
# rows 105..163 down to 106..164 and copies the row-105 formatting
# (incl. the date style on column D) onto the new, still-empty row.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new observation.
$ws.Range("A105").Value = 4
$ws.Range("B105").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C105").Value = "Los Lagos"
$ws.Range("D105").Value = 45089
$ws.Range("E105").Value = 10
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100104
$ws.Range("H105").Value = "Frutos de pepita"
$ws.Range("I105").Value = 100104003
$ws.Range("J105").Value = "Membrillo"
$ws.Range("K105").Value = "Champion"
$ws.Range("L105").Value = "Primera"
$ws.Range("M105").Value = 200
$ws.Range("N105").Value = 12000
$ws.Range("O105").Value = 13000
$ws.Range("P105").Value = 12500
$ws.Range("Q105").Value = "$/caja 18 kilos empedrada"
$ws.Range("R105").Value = "Región de O'Higgins"
$ws.Range("S105").Value = 694
$ws.Range("T105").Value = 18
